$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-04-07"
$ws.Range("B2").Value = "Monday"
$ws.Range("C2").Value = "Meenakshi"
$ws.Range("D2").Value = "Abhi pratap singh"
$ws.Range("E2").Value = "P"
